# Update the TV price-list sheet: rows 25-72 (one "page" of the list)
# are refreshed with the listing data that currently lives 48 rows below
# them (rows 73-120), mirroring the same "catalog rotation" already
# visible later in the sheet. This covers the product name (col A),
# the regular price (col B) and, where applicable, the discounted
# price / percentage (col C) and strike-through price (col D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$src = $ws.Range("A73:D120")
$dst = $ws.Range("A25:D72")

$dst.Value = $src.Value2
